# Update cryptocurrency price/volume data in the "cryptos" worksheet.
# Values that look like plain numbers (e.g. "1.002", "0.00001314") are
# prefixed with a leading apostrophe so Excel stores them as literal text
# (matching the source data, which keeps price strings like "23.889.16"
# or zero-padded values like "1.0000" that would otherwise be mangled by
# automatic number coercion).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.889.16"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.647.12"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'310.49"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.3892"
$ws.Range("D8").Value = "'0.3829"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").Value = "'51.11"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").Value = "'1.338"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").Value = "'1.001"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'0.08420"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "'23.81"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").Value = "'7.007"
$ws.Range("E14").Value = "  -3.36%  "
$ws.Range("D15").Value = "'7.853"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("D16").Value = "'0.00001314"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "1.650.19"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "'93.92"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "'0.06977"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "'19.52"
$ws.Range("E20").Value = "  -2.87%  "
$ws.Range("D21").Value = "'6.910"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "'13.64"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "23.887.86"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "'2.439"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("D26").Value = "'2.914"
$ws.Range("E26").Value = "  -7.95%  "
$ws.Range("D27").Value = "'21.91"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("D28").Value = "'153.72"
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("D29").Value = "'5.451"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").Value = "'136.85"
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("D31").Value = "'7.681"
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("D32").Value = "'2.507"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").Value = "1.834.04"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").Value = "'0.08106"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "'0.9840"
$ws.Range("E35").Value = "  -6.23%  "
$ws.Range("D36").Value = "'0.02918"
$ws.Range("E36").Value = "  -3.64%  "
$ws.Range("D37").Value = "'6.679"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").Value = "'0.2681"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").Value = "'10.45"
$ws.Range("E39").Value = "  -5.26%  "
$ws.Range("D40").Value = "'0.09091"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").Value = "'0.7527"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").Value = "'13.35"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").Value = "'1.422"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "'16.73"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("D45").Value = "'0.6917"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").Value = "'2.430"
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("D47").Value = "'4.094"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "'1.000"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").Value = "'0.08268"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").Value = "'134.57"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "'1.223"
$ws.Range("E51").Value = "  -1.04%  "
